# edit.ps1 - Applies the "Updated analyzer logic and config" changes
# to what_test_controls_analysis_results.xlsx
#
# Touches three worksheets:
#   - "Analysis Results"      (sheet1): column widths + recomputed scores
#   - "Keyword Matches"       (sheet2): shortened WHY-keyword excerpts
#   - "Enhancement Feedback"  (sheet3): column width + WHY feedback cleared to "None"
#   - "Executive Summary"     (sheet5): summary stats refreshed

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Analysis Results
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Analysis Results")

# Column width tweaks (G, J, K narrower/rebalanced)
$ws1.Columns.Item(7).ColumnWidth  = 20.0          # col 7  (G) -> raw width 20.9
$ws1.Columns.Item(10).ColumnWidth = 11.333333      # col 10 (J) -> raw width 12.1
$ws1.Columns.Item(11).ColumnWidth = 19.0          # col 11 (K) -> raw width 19.8

# Row 2 (CTRL-001)
$ws1.Range("G2").Value = 95
$ws1.Range("I2").Value = 70.59999999999999

# Row 3 (CTRL-002)
$ws1.Range("C3").Value = 12.4
$ws1.Range("E3").Value = "WHO, WHEN, WHAT"
$ws1.Range("I3").Value = 15
$ws1.Range("J3").Value = 63
$ws1.Range("K3").Value = 120

# Row 4 (CTRL-003)
$ws1.Range("C4").Value = 12.57341811320755
$ws1.Range("I4").Value = 57.5780603773585
$ws1.Range("J4").Value = 70
$ws1.Range("K4").Value = 10

# Row 5 (CTRL-004)
$ws1.Range("I5").Value = 43.86000000000001

# Row 6 (CTRL-005)
$ws1.Range("C6").Value = 61.175
$ws1.Range("G6").Value = 100
$ws1.Range("I6").Value = 92.25
$ws1.Range("J6").Value = 35

# Row 7 (CTRL-006)
$ws1.Range("G7").Value = 73.49999999999999
$ws1.Range("I7").Value = 75.00000000000001

# Row 8 (CTRL-007)
$ws1.Range("G8").Value = 100
$ws1.Range("H8").Value = 90
$ws1.Range("I8").Value = 86.41666666666666
$ws1.Range("J8").Value = 70

# Row 9 (CTRL-008)
$ws1.Range("I9").Value = 15

# Row 10 (CTRL-009)
$ws1.Range("G10").Value = 44.99999999999999
$ws1.Range("H10").Value = 90
$ws1.Range("I10").Value = 63.46774193548386
$ws1.Range("J10").Value = 90

# Row 11 (CTRL-010)
$ws1.Range("C11").Value = 21.84778
$ws1.Range("I11").Value = 49.4926
$ws1.Range("J11").Value = 90

# Row 12 (CTRL-011)
$ws1.Range("G12").Value = 95
$ws1.Range("H12").Value = 90
$ws1.Range("I12").Value = 45.08571428571428

# Row 13 (CTRL-012)
$ws1.Range("C13").Value = 43.24796
$ws1.Range("H13").Value = 85
$ws1.Range("I13").Value = 57.4932
$ws1.Range("J13").Value = 90

# Row 14 (CTRL-013)
$ws1.Range("C14").Value = 10.8
$ws1.Range("E14").Value = "WHO, WHEN, WHAT, ESCALATION"
$ws1.Range("I14").Value = 15
$ws1.Range("J14").Value = 63

# Row 15 (CTRL-014)
$ws1.Range("C15").Value = 11.5
$ws1.Range("I15").Value = 15
$ws1.Range("J15").Value = 70

# Row 16 (CTRL-015)
$ws1.Range("C16").Value = 24.0726
$ws1.Range("I16").Value = 50.242
$ws1.Range("J16").Value = 90

# ---------------------------------------------------------------------------
# Sheet: Keyword Matches
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Keyword Matches")

$ws2.Range("E3").Value  = "to ensure thresholds are set to appropriate limits and manag"
$ws2.Range("E6").Value  = "to ensure accuracy and completeness"
$ws2.Range("E10").Value = "to ensure all items are accounted for"
$ws2.Range("E11").Value = "to ensure timely removal of access"
$ws2.Range("E13").Value = "to ensure compliance with regulatory requirements"
$ws2.Range("E14").Value = "to ensure all changes are properly authorized"
$ws2.Range("E16").Value = "to ensure all exceptions are resolved before start of busine"

# ---------------------------------------------------------------------------
# Sheet: Enhancement Feedback
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Enhancement Feedback")

$ws3.Columns.Item(5).ColumnWidth = 14.5    # col 5 (E) -> raw width 15.4

for ($r = 2; $r -le 16; $r++) {
    $ws3.Range("E$r").Value = "None"
}

# ---------------------------------------------------------------------------
# Sheet: Executive Summary
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Executive Summary")

$ws5.Range("B4").Value = "'32.1"
$ws5.Range("B17").Value = "5 (33.3%)"
$ws5.Range("B18").Value = "14 (93.3%)"
